# Auto-generated Excel COM-interop script to apply the Rafflesia_Profits update.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled pricing refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value2 = 184
$ws.Range("I8").Value2 = 184
$ws.Range("K8").Value2 = 552
$ws.Range("M8").Value2 = -413
$ws.Range("H31").Value2 = 574.5
$ws.Range("I31").Value2 = 574.5
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 1723.5
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = -1493.5
$ws.Range("N31").ClearContents()
$ws.Range("H38").Value2 = 2329.4443
$ws.Range("I38").Value2 = 193.6
$ws.Range("K38").Value2 = 580.8
$ws.Range("M38").Value2 = -208.8
$ws.Range("H39").Value2 = 25.5
$ws.Range("I39").Value2 = 25.666666
$ws.Range("K39").Value2 = 76.99999800000001
$ws.Range("M39").Value2 = 219.000002
$ws.Range("H40").Value2 = 1880
$ws.Range("I40").Value2 = 2000
$ws.Range("K40").Value2 = 2000
$ws.Range("M40").Value2 = -1825
$ws.Range("H42").Value2 = 1811
$ws.Range("I42").Value2 = 100
$ws.Range("J42").Value2 = 2666.5
$ws.Range("K42").Value2 = 300
$ws.Range("L42").Value2 = 7999.5
$ws.Range("M42").Value2 = -70
$ws.Range("N42").Value2 = -8459.5
$ws.Range("H45").Value2 = 999
$ws.Range("I45").Value2 = 999
$ws.Range("J45").Value2 = 0
$ws.Range("K45").Value2 = 2997
$ws.Range("L45").Value2 = 0
$ws.Range("M45").Value2 = -2805
$ws.Range("N45").ClearContents()
$ws.Range("H80").Value2 = 1829.5454
$ws.Range("J80").Value2 = 2230.1667
$ws.Range("L80").Value2 = 6690.500100000001
$ws.Range("N80").Value2 = -8686.500100000001
$ws.Range("H82").Value2 = 0
$ws.Range("I82").Value2 = 0
$ws.Range("K82").Value2 = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value2 = 1829.5454
$ws.Range("J83").Value2 = 2230.1667
$ws.Range("L83").Value2 = 20071.5003
$ws.Range("N83").Value2 = -30055.5003
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("M85").ClearContents()
$ws.Range("H137").Value2 = 1632.3334
$ws.Range("I137").Value2 = 1632.3334
$ws.Range("K137").Value2 = 4897.0002
$ws.Range("M137").Value2 = -2347.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1358.3334
$ws.Range("I2").Value2 = 1358.3334
$ws.Range("K2").Value2 = 1358.3334
$ws.Range("M2").Value2 = -1245.3334
$ws.Range("H61").Value2 = 2968.7273
$ws.Range("I61").Value2 = 2968.7273
$ws.Range("K61").Value2 = 2968.7273
$ws.Range("M61").Value2 = -2756.7273
$ws.Range("H116").Value2 = 1358.3334
$ws.Range("I116").Value2 = 1358.3334
$ws.Range("K116").Value2 = 1358.3334
$ws.Range("M116").Value2 = 935.6666
$ws.Range("H132").Value2 = 3079.7
$ws.Range("I132").Value2 = 2310.7778
$ws.Range("K132").Value2 = 6932.3334
$ws.Range("M132").Value2 = -4402.3334
$ws.Range("H136").Value2 = 2968.7273
$ws.Range("I136").Value2 = 2968.7273
$ws.Range("K136").Value2 = 8906.1819
$ws.Range("M136").Value2 = -6356.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1358.3334
$ws.Range("I3").Value2 = 1358.3334
$ws.Range("K3").Value2 = 1358.3334
$ws.Range("M3").Value2 = -1244.3334
$ws.Range("H29").Value2 = 19499.5
$ws.Range("J29").Value2 = 9999
$ws.Range("L29").Value2 = 9999
$ws.Range("N29").Value2 = -10577
$ws.Range("H36").Value2 = 10666
$ws.Range("I36").Value2 = 9999
$ws.Range("J36").Value2 = 12000
$ws.Range("K36").Value2 = 9999
$ws.Range("L36").Value2 = 12000
$ws.Range("M36").Value2 = -9465
$ws.Range("N36").Value2 = -13068
$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 0
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 0
$ws.Range("L89").Value2 = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value2 = 2132
$ws.Range("I134").Value2 = 2118.3333
$ws.Range("J134").Value2 = 2214
$ws.Range("K134").Value2 = 6354.999899999999
$ws.Range("L134").Value2 = 6642
$ws.Range("M134").Value2 = -3819.999899999999
$ws.Range("N134").Value2 = -11712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value2 = 35499.5
$ws.Range("I25").Value2 = 35499.5
$ws.Range("K25").Value2 = 35499.5
$ws.Range("M25").Value2 = -35325.5
$ws.Range("H86").Value2 = 6004
$ws.Range("I86").Value2 = 5000
$ws.Range("K86").Value2 = 5000
$ws.Range("M86").Value2 = -3877
$ws.Range("H89").Value2 = 6004
$ws.Range("I89").Value2 = 5000
$ws.Range("K89").Value2 = 25000
$ws.Range("M89").Value2 = -19384
$ws.Range("H99").Value2 = 5210.278
$ws.Range("I99").Value2 = 4739.8184
$ws.Range("K99").Value2 = 4739.8184
$ws.Range("M99").Value2 = -3241.8184
$ws.Range("H107").Value2 = 781.9231
$ws.Range("I107").Value2 = 791.2222
$ws.Range("K107").Value2 = 791.2222
$ws.Range("M107").Value2 = 1128.7778
$ws.Range("H126").Value2 = 5210.278
$ws.Range("I126").Value2 = 4739.8184
$ws.Range("K126").Value2 = 14219.4552
$ws.Range("M126").Value2 = -11749.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value2 = 700
$ws.Range("I98").Value2 = 700
$ws.Range("K98").Value2 = 2100
$ws.Range("M98").Value2 = -602
$ws.Range("H107").Value2 = 2612.111
$ws.Range("J107").Value2 = 2900.5
$ws.Range("L107").Value2 = 8701.5
$ws.Range("N107").Value2 = -12541.5
$ws.Range("H131").Value2 = 2574.9375
$ws.Range("I131").Value2 = 250
$ws.Range("J131").Value2 = 2907.0715
$ws.Range("K131").Value2 = 750
$ws.Range("L131").Value2 = 8721.2145
$ws.Range("M131").Value2 = 4290
$ws.Range("N131").Value2 = -18801.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 18215
$ws.Range("I80").Value2 = 2372.5
$ws.Range("K80").Value2 = 2372.5
$ws.Range("M80").Value2 = -1374.5
$ws.Range("H83").Value2 = 18215
$ws.Range("I83").Value2 = 2372.5
$ws.Range("K83").Value2 = 11862.5
$ws.Range("M83").Value2 = -6870.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 4711.7646
$ws.Range("I16").Value2 = 3531.818
$ws.Range("J16").Value2 = 6875
$ws.Range("K16").Value2 = 3531.818
$ws.Range("L16").Value2 = 6875
$ws.Range("M16").Value2 = -3361.818
$ws.Range("N16").Value2 = -7215
$ws.Range("H46").Value2 = 2000
$ws.Range("I46").Value2 = 2000
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 2000
$ws.Range("L46").Value2 = 0
$ws.Range("M46").Value2 = -1812
$ws.Range("N46").ClearContents()
$ws.Range("H68").Value2 = 0
$ws.Range("I68").Value2 = 0
$ws.Range("K68").Value2 = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value2 = 0
$ws.Range("I71").Value2 = 0
$ws.Range("K71").Value2 = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value2 = 1398.625
$ws.Range("I82").Value2 = 1398.625
$ws.Range("J82").Value2 = 0
$ws.Range("K82").Value2 = 1398.625
$ws.Range("L82").Value2 = 0
$ws.Range("M82").Value2 = -1037.625
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value2 = 1398.625
$ws.Range("I85").Value2 = 1398.625
$ws.Range("J85").Value2 = 0
$ws.Range("K85").Value2 = 1398.625
$ws.Range("L85").Value2 = 0
$ws.Range("M85").Value2 = -150.625
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value2 = 0
$ws.Range("I132").Value2 = 0
$ws.Range("K132").Value2 = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 2987.5
$ws.Range("I62").Value2 = 2987.5
$ws.Range("K62").Value2 = 2987.5
$ws.Range("M62").Value2 = -2363.5
$ws.Range("H65").Value2 = 2987.5
$ws.Range("I65").Value2 = 2987.5
$ws.Range("K65").Value2 = 14937.5
$ws.Range("M65").Value2 = -11817.5
